$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2587124234681539
$ws.Range("C2").Value = 0.03611214932875839
$ws.Range("E2").Value = 0.4900217385315528
$ws.Range("F2").Value = 2.229381280129004
$ws.Range("G2").Value = 0.002448025733672279
$ws.Range("I2").Value = 0.6208505565256885
$ws.Range("J2").Value = 0.05585128492409552
$ws.Range("K2").Value = 0.2642209713797286
$ws.Range("M2").Value = 0.4021137604249532
$ws.Range("O2").Value = 2.779098296046641

$ws.Range("B3").Value = 0.227432949004708
$ws.Range("C3").Value = 0.03186824337062433
$ws.Range("E3").Value = 0.4800693629255974
$ws.Range("F3").Value = 2.221081260319067
$ws.Range("G3").Value = 0.0024502086728924
$ws.Range("I3").Value = 0.6285404289231629
$ws.Range("J3").Value = 0.05622708831677414
$ws.Range("K3").Value = 0.231054985630891
$ws.Range("M3").Value = 0.3801312215460229
$ws.Range("O3").Value = 2.806929025359267

$ws.Range("B4").Value = 0.2081939692552055
$ws.Range("C4").Value = 0.02924833672945226
$ws.Range("E4").Value = 0.4741820905040441
$ws.Range("F4").Value = 2.217216805839513
$ws.Range("G4").Value = 0.002451619739001172
$ws.Range("I4").Value = 0.6336030001194057
$ws.Range("J4").Value = 0.05647287603985962
$ws.Range("K4").Value = 0.2106370737693339
$ws.Range("M4").Value = 0.3667744712635468
$ws.Range("O4").Value = 2.825623226152146

$ws.Range("B5").Value = 0.20034605430223
$ws.Range("C5").Value = 0.02817719804622243
$ws.Range("E5").Value = 0.4718392699871998
$ws.Range("F5").Value = 2.215951769323766
$ws.Range("G5").Value = 0.002452212598247469
$ws.Range("I5").Value = 0.6357517339845113
$ws.Range("J5").Value = 0.05657682577883882
$ws.Range("K5").Value = 0.2023035496410017
$ws.Range("M5").Value = 0.3613670992268609
$ws.Range("O5").Value = 2.833644900320891

$ws.Range("B6").Value = 0.1990424528122219
$ws.Range("C6").Value = 0.0279991259816228
$ws.Range("E6").Value = 0.4714536487231058
$ws.Range("F6").Value = 2.215760420257155
$ws.Range("G6").Value = 0.002452312121041236
$ws.Range("I6").Value = 0.6361137041510379
$ws.Range("J6").Value = 0.05659431560868899
$ws.Range("K6").Value = 0.2009189998841094
$ws.Range("M6").Value = 0.3604713675732611
$ws.Range("O6").Value = 2.835001265622097

$ws.Range("B7").Value = 0.208088160764504
$ws.Range("C7").Value = 0.02923390510825641
$ws.Range("E7").Value = 0.4741502663774071
$ws.Range("F7").Value = 2.217198490889032
$ws.Range("G7").Value = 0.002451627662069023
$ws.Range("I7").Value = 0.6336316317759838
$ws.Range("J7").Value = 0.05647426259139632
$ws.Range("K7").Value = 0.2105247371805063
$ws.Range("M7").Value = 0.3667014010251819
$ws.Range("O7").Value = 2.825729775234521

$ws.Range("B8").Value = 0.2479344455516923
$ws.Range("C8").Value = 0.03465181093180547
$ws.Range("E8").Value = 0.4865438121792991
$ws.Range("F8").Value = 2.226263775718749
$ws.Range("G8").Value = 0.002448763762230688
$ws.Range("I8").Value = 0.623431234038069
$ws.Range("J8").Value = 0.05597774373788234
$ws.Range("K8").Value = 0.2527967989292108
$ws.Range("M8").Value = 0.3945051398658563
$ws.Range("O8").Value = 2.788360875155945

$ws.Range("B9").Value = 0.3257913842982703
$ws.Range("C9").Value = 0.04516257514771382
$ws.Range("E9").Value = 0.5126192294069938
$ws.Range("F9").Value = 2.253817866297126
$ws.Range("G9").Value = 0.002443706500191748
$ws.Range("I9").Value = 0.6061353773584912
$ws.Range("J9").Value = 0.0551231328710351
$ws.Range("K9").Value = 0.335247500835095
$ws.Range("M9").Value = 0.4501358706769878
$ws.Range("O9").Value = 2.727833490151227

$ws.Range("B10").Value = 0.3828031701116572
$ws.Range("C10").Value = 0.05281407443531805
$ws.Range("E10").Value = 0.5328570540435322
$ws.Range("F10").Value = 2.2800328607631
$ws.Range("G10").Value = 0.002440328248076357
$ws.Range("I10").Value = 0.5950806302202736
$ws.Range("J10").Value = 0.05456741999789116
$ws.Range("K10").Value = 0.3955355662577631
$ws.Range("M10").Value = 0.4916767251135141
$ws.Range("O10").Value = 2.691152677413854

$ws.Range("B11").Value = 0.4086946458871807
$ws.Range("C11").Value = 0.05627932189648277
$ws.Range("E11").Value = 0.5422985258098407
$ws.Range("F11").Value = 2.293257965729268
$ws.Range("G11").Value = 0.002438863923373891
$ws.Range("I11").Value = 0.5904108438965565
$ws.Range("J11").Value = 0.05433019578169329
$ws.Range("K11").Value = 0.4228962328456873
$ws.Range("M11").Value = 0.5107189387872779
$ws.Range("O11").Value = 2.676159982997348

$ws.Range("B12").Value = 0.4184924084242709
$ws.Range("C12").Value = 0.05758926047641921
$ws.Range("E12").Value = 0.5459075440473384
$ws.Range("F12").Value = 2.2984529769586
$ws.Range("G12").Value = 0.002438319787202485
$ws.Range("I12").Value = 0.5886942011338938
$ws.Range("J12").Value = 0.05424259795329078
$ws.Range("K12").Value = 0.433247297900067
$ws.Range("M12").Value = 0.5179504007812596
$ws.Range("O12").Value = 2.670726445104947

$ws.Range("B13").Value = 0.416382594739531
$ws.Range("C13").Value = 0.05730724380028107
$ws.Range("E13").Value = 0.5451287781922645
$ws.Range("F13").Value = 2.297325822924989
$ws.Range("G13").Value = 0.002438436516044211
$ws.Range("I13").Value = 0.5890616102101021
$ws.Range("J13").Value = 0.05426136444851082
$ws.Range("K13").Value = 0.4310184561671804
$ws.Range("M13").Value = 0.5163920646807938
$ws.Range("O13").Value = 2.67188580320898

$ws.Range("B14").Value = 0.4095008525434878
$ws.Range("C14").Value = 0.05638713723547539
$ws.Range("E14").Value = 0.5425947662220381
$ws.Range("F14").Value = 2.29368161539486
$ws.Range("G14").Value = 0.00243881894925597
$ws.Range("I14").Value = 0.5902685780397299
$ws.Range("J14").Value = 0.05432294431712137
$ws.Range("K14").Value = 0.4237480211162108
$ws.Range("M14").Value = 0.5113134644410948
$ws.Range("O14").Value = 2.675708073134658

$ws.Range("B15").Value = 0.4052846885429631
$ws.Range("C15").Value = 0.05582324751314616
$ws.Range("E15").Value = 0.5410470019669731
$ws.Range("F15").Value = 2.29147378093603
$ws.Range("G15").Value = 0.002439054549956315
$ws.Range("I15").Value = 0.5910146161365155
$ws.Range("J15").Value = 0.05436095451532985
$ws.Range("K15").Value = 0.4192933745928826
$ws.Range("M15").Value = 0.5082053483758457
$ws.Range("O15").Value = 2.678081092814068

$ws.Range("B16").Value = 0.3811101814798974
$ws.Range("C16").Value = 0.0525872968092358
$ws.Range("E16").Value = 0.5322447571992512
$ws.Range("F16").Value = 2.279194727675915
$ws.Range("G16").Value = 0.0024404253992113
$ws.Range("I16").Value = 0.5953930431676469
$ws.Range("J16").Value = 0.05458323604613824
$ws.Range("K16").Value = 0.3937461359424788
$ws.Range("M16").Value = 0.4904351675932404
$ws.Range("O16").Value = 2.692166593175315

$ws.Range("B17").Value = 0.3662683951093868
$ws.Range("C17").Value = 0.05059814918479333
$ws.Range("E17").Value = 0.5269050434013565
$ws.Range("F17").Value = 2.271994861474369
$ws.Range("G17").Value = 0.00244128489515931
$ws.Range("I17").Value = 0.5981710832138987
$ws.Range("J17").Value = 0.05472358311961223
$ws.Range("K17").Value = 0.3780567928211838
$ws.Range("M17").Value = 0.4795706996147899
$ws.Range("O17").Value = 2.701241598336594

$ws.Range("B18").Value = 0.3577277316993843
$ws.Range("C18").Value = 0.04945259223205767
$ws.Range("E18").Value = 0.523855922721026
$ws.Range("F18").Value = 2.267976026906609
$ws.Range("G18").Value = 0.00244178607767593
$ws.Range("I18").Value = 0.59980273190844
$ws.Range("J18").Value = 0.05480577310923707
$ws.Range("K18").Value = 0.3690266411788627
$ws.Range("M18").Value = 0.4733354159544092
$ws.Range("O18").Value = 2.706620690159085

$ws.Range("B19").Value = 0.3548353285077326
$ws.Range("C19").Value = 0.04906447858131457
$ws.Range("E19").Value = 0.5228273481509973
$ws.Range("F19").Value = 2.266636328871712
$ws.Range("G19").Value = 0.00244195694269641
$ws.Range("I19").Value = 0.6003609820745837
$ws.Range("J19").Value = 0.05483385319310941
$ws.Range("K19").Value = 0.3659681628325586
$ws.Range("M19").Value = 0.4712266123425763
$ws.Range("O19").Value = 2.70846932325189

$ws.Range("B20").Value = 0.36784875266531
$ws.Range("C20").Value = 0.05081004798823585
$ws.Range("E20").Value = 0.5274711740823363
$ws.Range("F20").Value = 2.272748637662431
$ws.Range("G20").Value = 0.002441192694677774
$ws.Range("I20").Value = 0.5978718580991931
$ws.Range("J20").Value = 0.05470849124137622
$ws.Range("K20").Value = 0.3797275808113056
$ws.Range("M20").Value = 0.4807258278571425
$ws.Range("O20").Value = 2.700259049970683

$ws.Range("B21").Value = 0.4115223751619226
$ws.Range("C21").Value = 0.05665745700866864
$ws.Range("E21").Value = 0.5433381519764993
$ws.Range("F21").Value = 2.294746933327303
$ws.Range("G21").Value = 0.00243870633825354
$ws.Range("I21").Value = 0.589912658544641
$ws.Range("J21").Value = 0.05430479623447404
$ws.Range("K21").Value = 0.4258837949648182
$ws.Range("M21").Value = 0.5128046158260773
$ws.Range("O21").Value = 2.674578758581177

$ws.Range("B22").Value = 0.4400258655811058
$ws.Range("C22").Value = 0.06046578972657812
$ws.Range("E22").Value = 0.5539047513130271
$ws.Range("F22").Value = 2.310213805553488
$ws.Range("G22").Value = 0.002437141794422839
$ws.Range("I22").Value = 0.5850122642005253
$ws.Range("J22").Value = 0.05405397622158503
$ws.Range("K22").Value = 0.4559920361178911
$ws.Range("M22").Value = 0.5338898404150001
$ws.Range("O22").Value = 2.659216819864696

$ws.Range("B23").Value = 0.4248168264028607
$ws.Range("C23").Value = 0.05843444394393771
$ws.Range("E23").Value = 0.5482471951923884
$ws.Range("F23").Value = 2.301859120937266
$ws.Range("G23").Value = 0.002437971306771622
$ws.Range("I23").Value = 0.5876000955010845
$ws.Range("J23").Value = 0.05418665421827562
$ws.Range("K23").Value = 0.4399281365104457
$ws.Range("M23").Value = 0.5226253797496128
$ws.Range("O23").Value = 2.667285587504651

$ws.Range("B24").Value = 0.3671342976534788
$ws.Range("C24").Value = 0.05071425467110657
$ws.Range("E24").Value = 0.5272151617552225
$ws.Range("F24").Value = 2.272407480209807
$ws.Range("G24").Value = 0.002441234356627631
$ws.Range("I24").Value = 0.5980070302610372
$ws.Range("J24").Value = 0.05471530959867543
$ws.Range("K24").Value = 0.3789722491049758
$ws.Range("M24").Value = 0.4802035605822539
$ws.Range("O24").Value = 2.700702756263823

$ws.Range("B25").Value = 0.3047610365135256
$ws.Range("C25").Value = 0.04233145940700922
$ws.Range("E25").Value = 0.5053754515482751
$ws.Range("F25").Value = 2.245315951283331
$ws.Range("G25").Value = 0.002445015145287977
$ws.Range("I25").Value = 0.6105243249478782
$ws.Range("J25").Value = 0.05534162441227153
$ws.Range("K25").Value = 0.3129918136234835
$ws.Range("M25").Value = 0.4349683541639635
$ws.Range("O25").Value = 2.742841003481331

